$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.374193943068803
$ws.Cells.Item(2, 3).Value = 0.2344924882585957
$ws.Cells.Item(2, 5).Value = 0.2279316778279501
$ws.Cells.Item(2, 6).Value = 1.925295705760774
$ws.Cells.Item(2, 7).Value = 0.002447603289630147
$ws.Cells.Item(2, 10).Value = 0.03906683956150303
$ws.Cells.Item(2, 12).Value = 0.4946125413588902
$ws.Cells.Item(2, 14).Value = 1.309358077573656
$ws.Cells.Item(2, 15).Value = 3.200513877259709
$ws.Cells.Item(3, 2).Value = 1.275719050755981
$ws.Cells.Item(3, 3).Value = 0.2268029899658757
$ws.Cells.Item(3, 5).Value = 0.2283203109745315
$ws.Cells.Item(3, 6).Value = 1.921500310126433
$ws.Cells.Item(3, 7).Value = 0.002450545400425724
$ws.Cells.Item(3, 10).Value = 0.03867317094348266
$ws.Cells.Item(3, 12).Value = 0.4837710350266491
$ws.Cells.Item(3, 14).Value = 1.321086245987644
$ws.Cells.Item(3, 15).Value = 3.205088444612898
$ws.Cells.Item(4, 2).Value = 1.215578710516013
$ws.Cells.Item(4, 3).Value = 0.2220431398109497
$ws.Cells.Item(4, 5).Value = 0.2286342329531657
$ws.Cells.Item(4, 6).Value = 1.920228432977993
$ws.Cells.Item(4, 7).Value = 0.002452449050775297
$ws.Cells.Item(4, 10).Value = 0.03843551413975987
$ws.Cells.Item(4, 12).Value = 0.4772955714949916
$ws.Cells.Item(4, 14).Value = 1.328806723772871
$ws.Cells.Item(4, 15).Value = 3.209950593166553
$ws.Cells.Item(5, 2).Value = 1.191153878413161
$ws.Cells.Item(5, 3).Value = 0.2200938981304432
$ws.Cells.Item(5, 5).Value = 0.2287811331394778
$ws.Cells.Item(5, 6).Value = 1.919976348500967
$ws.Cells.Item(5, 7).Value = 0.002453249314695333
$ws.Cells.Item(5, 10).Value = 0.03833970188750335
$ws.Cells.Item(5, 12).Value = 0.474702578844159
$ws.Cells.Item(5, 14).Value = 1.332083554277077
$ws.Cells.Item(5, 15).Value = 3.212447775773143
$ws.Cells.Item(6, 2).Value = 1.187103200432261
$ws.Cells.Item(6, 3).Value = 0.2197696535966429
$ws.Cells.Item(6, 5).Value = 0.2288066729819018
$ws.Cells.Item(6, 6).Value = 1.919950570202474
$ws.Cells.Item(6, 7).Value = 0.00245338368053224
$ws.Cells.Item(6, 10).Value = 0.03832385526565218
$ws.Cells.Item(6, 12).Value = 0.4742747880812601
$ws.Cells.Item(6, 14).Value = 1.332635564326701
$ws.Cells.Item(6, 15).Value = 3.212893570207314
$ws.Cells.Item(7, 2).Value = 1.215248971313656
$ws.Cells.Item(7, 3).Value = 0.2220168902359063
$ws.Cells.Item(7, 5).Value = 0.2286361372206009
$ws.Cells.Item(7, 6).Value = 1.920223955323578
$ws.Cells.Item(7, 7).Value = 0.002452459743897212
$ws.Cells.Item(7, 10).Value = 0.03843421777327904
$ws.Cells.Item(7, 12).Value = 0.4772604156937774
$ws.Cells.Item(7, 14).Value = 1.328850387054445
$ws.Cells.Item(7, 15).Value = 3.209982183205796
$ws.Cells.Item(8, 2).Value = 1.340173472500965
$ws.Cells.Item(8, 3).Value = 0.2318491960363644
$ws.Cells.Item(8, 5).Value = 0.2280500736069584
$ws.Cells.Item(8, 6).Value = 1.923767379350551
$ws.Cells.Item(8, 7).Value = 0.002448597603709564
$ws.Cells.Item(8, 10).Value = 0.03893026911963915
$ws.Cells.Item(8, 12).Value = 0.4908369023891339
$ws.Cells.Item(8, 14).Value = 1.313294210556933
$ws.Cells.Item(8, 15).Value = 3.201664698878091
$ws.Cells.Item(9, 2).Value = 1.587666642907607
$ws.Cells.Item(9, 3).Value = 0.2508211807686109
$ws.Cells.Item(9, 5).Value = 0.2274966589475795
$ws.Cells.Item(9, 6).Value = 1.939115968954127
$ws.Cells.Item(9, 7).Value = 0.002441791691098559
$ws.Cells.Item(9, 10).Value = 0.03993463277479492
$ws.Cells.Item(9, 12).Value = 0.5188902271892886
$ws.Cells.Item(9, 14).Value = 1.286905745730003
$ws.Cells.Item(9, 15).Value = 3.201673794021588
$ws.Cells.Item(10, 2).Value = 1.770984083217741
$ws.Cells.Item(10, 3).Value = 0.2645674429502947
$ws.Cells.Item(10, 5).Value = 0.2274513358704731
$ws.Cells.Item(10, 6).Value = 1.955519525856488
$ws.Cells.Item(10, 7).Value = 0.0024372547005415
$ws.Cells.Item(10, 10).Value = 0.04069110462112846
$ws.Cells.Item(10, 12).Value = 0.5403640168591579
$ws.Cells.Item(10, 14).Value = 1.270023404549868
$ws.Cells.Item(10, 15).Value = 3.211673120413508
$ws.Cells.Item(11, 2).Value = 1.854692832317539
$ws.Cells.Item(11, 3).Value = 0.2707784187329025
$ws.Cells.Item(11, 5).Value = 0.2275087593019265
$ws.Cells.Item(11, 6).Value = 1.964096818834051
$ws.Cells.Item(11, 7).Value = 0.002435290306701383
$ws.Cells.Item(11, 10).Value = 0.04103913060181341
$ws.Cells.Item(11, 12).Value = 0.5503187200205701
$ws.Cells.Item(11, 14).Value = 1.262886306514609
$ws.Cells.Item(11, 15).Value = 3.218401712096437
$ws.Cells.Item(12, 2).Value = 1.886435491417899
$ws.Cells.Item(12, 3).Value = 0.2731241829237945
$ws.Cells.Item(12, 5).Value = 0.2275416897257401
$ws.Cells.Item(12, 6).Value = 1.967505231927845
$ws.Cells.Item(12, 7).Value = 0.002434560674799748
$ws.Cells.Item(12, 10).Value = 0.04117146646019521
$ws.Cells.Item(12, 12).Value = 0.5541148818247166
$ws.Cells.Item(12, 14).Value = 1.260261672224956
$ws.Cells.Item(12, 15).Value = 3.221263807092555
$ws.Cells.Item(13, 2).Value = 1.879597211127361
$ws.Cells.Item(13, 3).Value = 0.2726192581675093
$ws.Cells.Item(13, 5).Value = 0.2275341006802769
$ws.Cells.Item(13, 6).Value = 1.966764036165785
$ws.Cells.Item(13, 7).Value = 0.002434717181733286
$ws.Cells.Item(13, 10).Value = 0.04114294152621767
$ws.Cells.Item(13, 12).Value = 0.5532961342531308
$ws.Cells.Item(13, 14).Value = 1.260823464027652
$ws.Cells.Item(13, 15).Value = 3.220633422431092
$ws.Cells.Item(14, 2).Value = 1.857303447701781
$ws.Cells.Item(14, 3).Value = 0.2709715310563467
$ws.Cells.Item(14, 5).Value = 0.2275112445456848
$ws.Cells.Item(14, 6).Value = 1.964374016756466
$ws.Cells.Item(14, 7).Value = 0.00243522999438585
$ws.Cells.Item(14, 10).Value = 0.04105000709778395
$ws.Cells.Item(14, 12).Value = 0.5506305021684028
$ws.Cells.Item(14, 14).Value = 1.262668812527735
$ws.Cells.Item(14, 15).Value = 3.218630879071753
$ws.Cells.Item(15, 2).Value = 1.843653553817887
$ws.Cells.Item(15, 3).Value = 0.2699614399818131
$ws.Cells.Item(15, 5).Value = 0.2274987000844568
$ws.Cells.Item(15, 6).Value = 1.962930947130431
$ws.Cells.Item(15, 7).Value = 0.002435545959557456
$ws.Cells.Item(15, 10).Value = 0.0409931527252354
$ws.Cells.Item(15, 12).Value = 0.5490011736794997
$ws.Cells.Item(15, 14).Value = 1.263809303082617
$ws.Cells.Item(15, 15).Value = 3.217445191279154
$ws.Cells.Item(16, 2).Value = 1.765519758430059
$ws.Cells.Item(16, 3).Value = 0.2641606803128411
$ws.Cells.Item(16, 5).Value = 0.2274491503525518
$ws.Cells.Item(16, 6).Value = 1.954981421721186
$ws.Cells.Item(16, 7).Value = 0.002437385074851215
$ws.Cells.Item(16, 10).Value = 0.04066843753820137
$ws.Cells.Item(16, 12).Value = 0.5397171793396893
$ws.Cells.Item(16, 14).Value = 1.270500750748951
$ws.Cells.Item(16, 15).Value = 3.211277311984105
$ws.Cells.Item(17, 2).Value = 1.717667206714339
$ws.Cells.Item(17, 3).Value = 0.260591191143817
$ws.Cells.Item(17, 5).Value = 0.2274387151786605
$ws.Cells.Item(17, 6).Value = 1.950390289906309
$ws.Cells.Item(17, 7).Value = 0.002438538750143391
$ws.Cells.Item(17, 10).Value = 0.04047022408506251
$ws.Cells.Item(17, 12).Value = 0.5340692615286713
$ws.Cells.Item(17, 14).Value = 1.274744736854295
$ws.Cells.Item(17, 15).Value = 3.208052301515693
$ws.Cells.Item(18, 2).Value = 1.690173590065797
$ws.Cells.Item(18, 3).Value = 0.2585341440423718
$ws.Cells.Item(18, 5).Value = 0.2274400609364449
$ws.Cells.Item(18, 6).Value = 1.947854562747253
$ws.Cells.Item(18, 7).Value = 0.002439211683680052
$ws.Cells.Item(18, 10).Value = 0.04035658524366781
$ws.Cells.Item(18, 12).Value = 0.5308382659170263
$ws.Cells.Item(18, 14).Value = 1.277236853687334
$ws.Cells.Item(18, 15).Value = 3.206402498678102
$ws.Cells.Item(19, 2).Value = 1.680869908148054
$ws.Cells.Item(19, 3).Value = 0.257836984807426
$ws.Cells.Item(19, 5).Value = 0.2274417797006976
$ws.Cells.Item(19, 6).Value = 1.94701403759025
$ws.Cells.Item(19, 7).Value = 0.002439441138688239
$ws.Cells.Item(19, 10).Value = 0.04031817280611705
$ws.Cells.Item(19, 12).Value = 0.5297473256618872
$ws.Cells.Item(19, 14).Value = 1.278089416751328
$ws.Cells.Item(19, 15).Value = 3.205879116199952
$ws.Cells.Item(20, 2).Value = 1.722758108835251
$ws.Cells.Item(20, 3).Value = 0.2609715812473041
$ws.Cells.Item(20, 5).Value = 0.2274390658188103
$ws.Cells.Item(20, 6).Value = 1.950868159609854
$ws.Cells.Item(20, 7).Value = 0.002438414970316938
$ws.Cells.Item(20, 10).Value = 0.04049128623939424
$ws.Cells.Item(20, 12).Value = 0.5346686787075612
$ws.Cells.Item(20, 14).Value = 1.274287669994457
$ws.Cells.Item(20, 15).Value = 3.208374373567324
$ws.Cells.Item(21, 2).Value = 1.863850485534556
$ws.Cells.Item(21, 3).Value = 0.2714556775987944
$ws.Cells.Item(21, 5).Value = 0.2275176546701942
$ws.Cells.Item(21, 6).Value = 1.965071670695451
$ws.Cells.Item(21, 7).Value = 0.002435078982879178
$ws.Cells.Item(21, 10).Value = 0.04107728949625766
$ws.Cells.Item(21, 12).Value = 0.5514127447663668
$ws.Cells.Item(21, 14).Value = 1.262124671358933
$ws.Cells.Item(21, 15).Value = 3.219210544089378
$ws.Cells.Item(22, 2).Value = 1.956318089085642
$ws.Cells.Item(22, 3).Value = 0.2782714615616158
$ws.Cells.Item(22, 5).Value = 0.2276341972002491
$ws.Cells.Item(22, 6).Value = 1.975289298433125
$ws.Cells.Item(22, 7).Value = 0.002432981701367984
$ws.Cells.Item(22, 10).Value = 0.04146345177336741
$ws.Cells.Item(22, 12).Value = 0.562510496654923
$ws.Cells.Item(22, 14).Value = 1.254630261118628
$ws.Cells.Item(22, 15).Value = 3.22812386742612
$ws.Cells.Item(23, 2).Value = 1.906943521665426
$ws.Cells.Item(23, 3).Value = 0.2746370991034155
$ws.Cells.Item(23, 5).Value = 0.2275660443974949
$ws.Cells.Item(23, 6).Value = 1.969750419832621
$ws.Cells.Item(23, 7).Value = 0.00243409348970973
$ws.Cells.Item(23, 10).Value = 0.0412570642065333
$ws.Cells.Item(23, 12).Value = 0.5565733563601469
$ws.Cells.Item(23, 14).Value = 1.258588558456907
$ws.Cells.Item(23, 15).Value = 3.223198879280119
$ws.Cells.Item(24, 2).Value = 1.720456457500518
$ws.Cells.Item(24, 3).Value = 0.2607996221360622
$ws.Cells.Item(24, 5).Value = 0.227438884412873
$ws.Cells.Item(24, 6).Value = 1.950651791461368
$ws.Cells.Item(24, 7).Value = 0.002438470900974425
$ws.Cells.Item(24, 10).Value = 0.04048176305158435
$ws.Cells.Item(24, 12).Value = 0.5343976321360628
$ws.Cells.Item(24, 14).Value = 1.274494147339055
$ws.Cells.Item(24, 15).Value = 3.208228128454692
$ws.Cells.Item(25, 2).Value = 1.520448948197327
$ws.Cells.Item(25, 3).Value = 0.2457222667990067
$ws.Cells.Item(25, 5).Value = 0.2275827869302276
$ws.Cells.Item(25, 6).Value = 1.934063864955348
$ws.Cells.Item(25, 7).Value = 0.0024435511693418
$ws.Cells.Item(25, 10).Value = 0.03965961473021196
$ws.Cells.Item(25, 12).Value = 0.5111488804184745
$ws.Cells.Item(25, 14).Value = 1.293604257801377
$ws.Cells.Item(25, 15).Value = 3.19991968341742
